$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("2:2").Insert()
$ws.Rows("2:2").ClearFormats()

$ws.Range("A2").Value = 6
$ws.Range("B2").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C2").Value = "Metropolitana"
$ws.Range("D2").Value = 44643
$ws.Range("D2").NumberFormat = $ws.Range("D3").NumberFormat
$ws.Range("E2").Value = 13
$ws.Range("F2").Value = "Fruta"
$ws.Range("G2").Value = 100104
$ws.Range("H2").Value = "Frutos de pepita"
$ws.Range("I2").Value = 100104003
$ws.Range("J2").Value = "Membrillo"
$ws.Range("K2").Value = "Champion"
$ws.Range("L2").Value = "Especial"
$ws.Range("M2").Value = 8
$ws.Range("N2").Value = 280000
$ws.Range("O2").Value = 280000
$ws.Range("P2").Value = 280000
$ws.Range("Q2").Value = "`$/bins (450 kilos)"
$ws.Range("R2").Value = "Región de O'Higgins"
$ws.Range("S2").Value = 622
$ws.Range("T2").Value = 450
